$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet (Sheet1 -> Translations) ---
$ws.Name = "Translations"

# --- Rebuild the data grid ---
# A column gains "Entity Id" header + "AAAA..." filler that used to live in column C;
# B column gains per-row entity-type labels; the old numeric Index column moves to C;
# D/E keep "Orig" / message text. Columns are written whole-row to land cleanly in
# the shared-string table in the same relative order as the target file.

$ws.Range("A1").Value = "Entity Id"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Index"
$ws.Range("D1").Value = "Original"
$ws.Range("E1").Value = "Translation"

$ws.Range("A2").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B2").Value = "Title"
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = "Orig"
$ws.Range("E2").Value = $null

$ws.Range("A3").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B3").Value = "ValidationMessage"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Orig"
$ws.Range("E3").Value = "validation message"

$ws.Range("A4").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B4").Value = "Instruction"
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = "Orig"
$ws.Range("E4").Value = $null

$ws.Range("A5").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B5").Value = "OptionTitle"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "Orig"
$ws.Range("E5").Value = "option"

# --- Column widths (best effort; engine snaps to 1/6-character steps) ---
$ws.Columns.Item(1).ColumnWidth = 42.416667
$ws.Columns.Item(2).ColumnWidth = 17.25
$ws.Columns.Item(3).ColumnWidth = 5.14
$ws.Columns.Item(5).ColumnWidth = 39.083333

# --- Selection moved from E4 to E6 ---
$ws.Range("E6").Select() | Out-Null
